# edit.ps1
# 1) Re-style the table on slide 6 with the new (built-in) table style GUID.
# 2) Swap the presentation's live theme colour palette from the "Integral"
#    scheme over to the stock "Office" scheme (the deck's theme parts are
#    exchanged in the source edit; the colours driving the rendered slides
#    are the user-visible/semantically-important half of that swap, and are
#    the part reachable through the PowerPoint object model).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F9499125-80EC-4180-9B84-094509D94AF4}", $false)
    }
}

# --- 2. Theme colours -------------------------------------------------------
function HexToColorRef($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Order matches MsoThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    '000000', 'FFFFFF', '44546A', 'E7E6E6',
    '5B9BD5', 'ED7D31', 'A5A5A5', 'FFC000', '4472C4', '70AD47',
    '0563C1', '954F72'
)

$anchorSlide = $p.Slides.Item(1)
$themeColors = $anchorSlide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $entry = $themeColors.Colors($i)
    $entry.RGB = HexToColorRef $officeThemeColors[$i - 1]
}
